$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.476.82"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "1.825.18"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -1.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4587"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3819"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.89%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07869"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9604"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.08"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").Value = "1.838.81"
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.853"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.092"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06593"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001021"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").Value = "27.451.81"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.296"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.255"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.13%  "
$ws.Range("D26").Value = "2.055.16"
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.35"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.049"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.279"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09313"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9346"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.565"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.232"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.317"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05929"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02189"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.129"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.144"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5762"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1819"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.960"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5407"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.83"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.876"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.72%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.37%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06589"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("E51").Value = "  -33.60%  "
